$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2017.6
$ws.Range("J17").Value = 2043.6207
$ws.Range("L17").Value = 6130.8621
$ws.Range("N17").Value = -6466.8621

$ws.Range("H32").Value = 661.7273
$ws.Range("J32").Value = 750
$ws.Range("L32").Value = 750
$ws.Range("N32").Value = -1402

$ws.Range("H76").Value = 4855.5
$ws.Range("I76").Value = 4139.0835
$ws.Range("K76").Value = 4139.0835
$ws.Range("M76").Value = -3824.0835

$ws.Range("H79").Value = 4855.5
$ws.Range("I79").Value = 4139.0835
$ws.Range("K79").Value = 4139.0835
$ws.Range("M79").Value = -3047.0835

$ws.Range("H87").Value = 64266.047
$ws.Range("J87").Value = 69199.28
$ws.Range("L87").Value = 69199.28
$ws.Range("N87").Value = -71695.28

$ws.Range("H90").Value = 64266.047
$ws.Range("J90").Value = 69199.28
$ws.Range("L90").Value = 207597.84
$ws.Range("N90").Value = -220077.84

$ws.Range("H112").Value = 1321.6177
$ws.Range("J112").Value = 1349.2188
$ws.Range("L112").Value = 4047.6564
$ws.Range("N112").Value = -6263.6564

$ws.Range("H116").Value = 58963.69
$ws.Range("J116").Value = 21961.334
$ws.Range("L116").Value = 21961.334
$ws.Range("N116").Value = -28845.334

$ws.Range("H127").Value = 899
$ws.Range("I127").Value = 899
$ws.Range("K127").Value = 2697
$ws.Range("M127").Value = 2263

$ws.Range("H132").Value = 2442.5483
$ws.Range("I132").Value = 2275.8215
$ws.Range("K132").Value = 6827.4645
$ws.Range("M132").Value = -4297.4645

$ws.Range("H136").Value = 99330.336
$ws.Range("J136").Value = 99330.336
$ws.Range("L136").Value = 99330.336
$ws.Range("N136").Value = -109530.336

$ws.Range("H137").Value = 2604.7827
$ws.Range("I137").Value = 1582.3334
$ws.Range("K137").Value = 4747.0002
$ws.Range("M137").Value = -2197.0002

$ws.Range("H138").Value = 1945.7222
$ws.Range("I138").Value = 1088.7333
$ws.Range("K138").Value = 3266.199900000001
$ws.Range("M138").Value = 1873.800099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2895.7314
$ws.Range("I32").Value = 2645.3872
$ws.Range("J32").Value = 6000
$ws.Range("K32").Value = 2645.3872
$ws.Range("L32").Value = 6000
$ws.Range("M32").Value = -2358.3872
$ws.Range("N32").Value = -6574

$ws.Range("H45").Value = 2594.4783
$ws.Range("I45").Value = 2898.75
$ws.Range("J45").Value = 1899
$ws.Range("K45").Value = 2898.75
$ws.Range("L45").Value = 1899
$ws.Range("M45").Value = -2521.75
$ws.Range("N45").Value = -2653

$ws.Range("H64").Value = 71199.39999999999
$ws.Range("J64").Value = 72666
$ws.Range("L64").Value = 72666
$ws.Range("N64").Value = -73162

$ws.Range("H67").Value = 71199.39999999999
$ws.Range("J67").Value = 72666
$ws.Range("L67").Value = 72666
$ws.Range("N67").Value = -74382

$ws.Range("H124").Value = 40000
$ws.Range("J124").Value = 40000
$ws.Range("L124").Value = 40000
$ws.Range("N124").Value = -49820

$ws.Range("H125").Value = 55000
$ws.Range("J125").Value = 55000
$ws.Range("L125").Value = 55000
$ws.Range("N125").Value = -64840

$ws.Range("H139").Value = 89999.664
$ws.Range("J139").Value = 89999.664
$ws.Range("L139").Value = 89999.664
$ws.Range("N139").Value = -100279.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1793.091
$ws.Range("J105").Value = 1879
$ws.Range("L105").Value = 1879
$ws.Range("N105").Value = -5373

$ws.Range("H107").Value = 3923.9546
$ws.Range("J107").Value = 4033.25
$ws.Range("L107").Value = 4033.25
$ws.Range("N107").Value = -7873.25

$ws.Range("H140").Value = 154971.75
$ws.Range("J140").Value = 154971.75
$ws.Range("L140").Value = 154971.75
$ws.Range("N140").Value = -165331.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5082
$ws.Range("I31").Value = 2281.8948
$ws.Range("K31").Value = 2281.8948
$ws.Range("M31").Value = -1986.8948

$ws.Range("H34").Value = 5082
$ws.Range("I34").Value = 2281.8948
$ws.Range("K34").Value = 2281.8948
$ws.Range("M34").Value = -2079.8948

$ws.Range("H68").Value = 73396.55499999999
$ws.Range("J68").Value = 73396.55499999999
$ws.Range("L68").Value = 73396.55499999999
$ws.Range("N68").Value = -74894.55499999999

$ws.Range("H71").Value = 73396.55499999999
$ws.Range("J71").Value = 73396.55499999999
$ws.Range("L71").Value = 220189.665
$ws.Range("N71").Value = -227677.665

$ws.Range("H122").Value = 1876.65
$ws.Range("I122").Value = 1689.0857
$ws.Range("K122").Value = 5067.257100000001
$ws.Range("M122").Value = -2617.257100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3103
$ws.Range("I5").Value = 993
$ws.Range("K5").Value = 2979
$ws.Range("M5").Value = -2867

$ws.Range("H37").Value = 99998
$ws.Range("J37").Value = 99998
$ws.Range("L37").Value = 299994
$ws.Range("N37").Value = -300218

$ws.Range("H80").Value = 2999.5
$ws.Range("J80").Value = 3199.6
$ws.Range("L80").Value = 9598.799999999999
$ws.Range("N80").Value = -11470.8

$ws.Range("H83").Value = 2999.5
$ws.Range("J83").Value = 3199.6
$ws.Range("L83").Value = 28796.4
$ws.Range("N83").Value = -38156.39999999999

$ws.Range("H135").Value = 3103
$ws.Range("I135").Value = 993
$ws.Range("K135").Value = 8937
$ws.Range("M135").Value = -6402

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3737.2307
$ws.Range("I102").Value = 2655
$ws.Range("K102").Value = 2655
$ws.Range("M102").Value = -1033

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6293.6924
$ws.Range("I40").Value = 6357.5713
$ws.Range("J40").Value = 6219.1665
$ws.Range("K40").Value = 6357.5713
$ws.Range("L40").Value = 6219.1665
$ws.Range("M40").Value = -6221.5713
$ws.Range("N40").Value = -6491.1665

$ws.Range("H46").Value = 1525.9231
$ws.Range("I46").Value = 1165.6364
$ws.Range("J46").Value = 1790.1333
$ws.Range("K46").Value = 1165.6364
$ws.Range("L46").Value = 1790.1333
$ws.Range("M46").Value = -977.6364000000001
$ws.Range("N46").Value = -2166.1333

$ws.Range("H82").Value = 3278.7144
$ws.Range("I82").Value = 886.1429000000001
$ws.Range("J82").Value = 5671.2856
$ws.Range("K82").Value = 886.1429000000001
$ws.Range("L82").Value = 5671.2856
$ws.Range("M82").Value = -525.1429000000001
$ws.Range("N82").Value = -6393.2856

$ws.Range("H85").Value = 3278.7144
$ws.Range("I85").Value = 886.1429000000001
$ws.Range("J85").Value = 5671.2856
$ws.Range("K85").Value = 886.1429000000001
$ws.Range("L85").Value = 5671.2856
$ws.Range("M85").Value = 361.8570999999999
$ws.Range("N85").Value = -8167.2856

$ws.Range("H132").Value = 2704.5
$ws.Range("I132").Value = 2460.4468
$ws.Range("K132").Value = 7381.340400000001
$ws.Range("M132").Value = -4851.340400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 75879.39999999999
$ws.Range("J75").Value = 77643.78
$ws.Range("L75").Value = 77643.78
$ws.Range("N75").Value = -79515.78

$ws.Range("H78").Value = 75879.39999999999
$ws.Range("J78").Value = 77643.78
$ws.Range("L78").Value = 232931.34
$ws.Range("N78").Value = -242291.34

$ws.Range("H122").Value = 3011
$ws.Range("I122").Value = 2886.4211
$ws.Range("K122").Value = 8659.263300000001
$ws.Range("M122").Value = -6209.263300000001

$ws.Range("H126").Value = 3579
$ws.Range("I126").Value = 3726.6667
$ws.Range("K126").Value = 11180.0001
$ws.Range("M126").Value = -8710.000100000001

$ws.Range("H132").Value = 1892.0541
$ws.Range("I132").Value = 1420.7119
$ws.Range("K132").Value = 4262.1357
$ws.Range("M132").Value = -1732.1357

$ws.Range("H136").Value = 20320166
$ws.Range("I136").Value = 26570602
$ws.Range("K136").Value = 79711806
$ws.Range("M136").Value = -79709256
